$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (A1:G1)
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Copy the header style (bold, bordered, centered) from C1 to the new header cells
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:G1").PasteSpecial(-4122) | Out-Null

# Updated prediction table: Rank, Team, WIN, TOP4, TOP5, RELEGATION, ExpPoints
$teams = @(
    @{Row=2;  Team="Arsenal";                  ExpPoints=79.05744001808478},
    @{Row=3;  Team="Manchester City";           ExpPoints=73.00153604450263},
    @{Row=4;  Team="Liverpool";                 ExpPoints=69.54971321713377},
    @{Row=5;  Team="Chelsea";                   ExpPoints=60.86149552075103},
    @{Row=6;  Team="Aston Villa";                ExpPoints=58.3300199739688},
    @{Row=7;  Team="Crystal Palace";            ExpPoints=57.97584102571155},
    @{Row=8;  Team="Newcastle United";           ExpPoints=57.33954485698675},
    @{Row=9;  Team="Brighton & Hove Albion";     ExpPoints=53.82433934419745},
    @{Row=10; Team="Tottenham Hotspur";          ExpPoints=53.05358940826264},
    @{Row=11; Team="AFC Bournemouth";            ExpPoints=51.35575031880398},
    @{Row=12; Team="Manchester United";          ExpPoints=50.32723163588085},
    @{Row=13; Team="Brentford";                  ExpPoints=49.8128619518311},
    @{Row=14; Team="Everton";                    ExpPoints=44.34672107679917},
    @{Row=15; Team="Fulham";                     ExpPoints=43.70372068442524},
    @{Row=16; Team="Nottingham Forest";          ExpPoints=42.41929178661415},
    @{Row=17; Team="Sunderland";                 ExpPoints=39.1058609164503},
    @{Row=18; Team="West Ham United";            ExpPoints=37.26803431118558},
    @{Row=19; Team="Leeds United";               ExpPoints=36.52626142352296},
    @{Row=20; Team="Burnley";                    ExpPoints=34.46308046442598},
    @{Row=21; Team="Wolverhampton Wanderers";    ExpPoints=32.32431463220376}
)

foreach ($t in $teams) {
    $r = $t.Row
    $ws.Cells.Item($r, 2).Value = $t.Team
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = $t.ExpPoints
}

$wb.Save()
